# Applies:
#  1) updates to the "time_taken" column (F2:F30) on the "data" sheet
#  2) a new "metadata" worksheet describing the panel query

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1) Refresh F2:F30 "time_taken" timestamps -----------------------------
$times = @(
    "2021-10-05 14:35:09.107804",
    "2021-10-05 14:35:09.107812",
    "2021-10-05 14:35:09.107816",
    "2021-10-05 14:35:09.107819",
    "2021-10-05 14:35:09.107822",
    "2021-10-05 14:35:09.107825",
    "2021-10-05 14:35:09.107827",
    "2021-10-05 14:35:09.107830",
    "2021-10-05 14:35:09.107832",
    "2021-10-05 14:35:09.107835",
    "2021-10-05 14:35:09.107838",
    "2021-10-05 14:35:09.107840",
    "2021-10-05 14:35:09.107843",
    "2021-10-05 14:35:09.107845",
    "2021-10-05 14:35:09.107848",
    "2021-10-05 14:35:09.107850",
    "2021-10-05 14:35:09.107853",
    "2021-10-05 14:35:09.107856",
    "2021-10-05 14:35:09.107858",
    "2021-10-05 14:35:09.107861",
    "2021-10-05 14:35:09.107864",
    "2021-10-05 14:35:09.107866",
    "2021-10-05 14:35:09.107869",
    "2021-10-05 14:35:09.107871",
    "2021-10-05 14:35:09.107874",
    "2021-10-05 14:35:09.107877",
    "2021-10-05 14:35:09.107879",
    "2021-10-05 14:35:09.107882",
    "2021-10-05 14:35:09.107884"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $times[$i]
}

# --- 2) Add a "metadata" worksheet after "data" -----------------------------
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"
$meta.Range("B1:G1").Font.Bold = $true

$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 1).Font.Bold = $true
$meta.Cells.Item(2, 2).Value = "Pain syndromes"
$meta.Cells.Item(2, 3).Value = 3126
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "0.28"
$meta.Cells.Item(2, 5).Value = "2021-06-11T05:55:44.605450Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:35:09.104089"
$meta.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/3126/?format=json"
